$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 1.022287785606226
$ws.Cells.Item(2, 4).Value = 1.017212861301559
$ws.Cells.Item(2, 5).Value = 1.008764243556
$ws.Cells.Item(2, 6).Value = 1.008320882087645
$ws.Cells.Item(2, 7).Value = 1.008018130962079
$ws.Cells.Item(2, 8).Value = 1.007668345802915
$ws.Cells.Item(2, 9).Value = 1.007425470330303
$ws.Cells.Item(2, 10).Value = 1.007462419686843
$ws.Cells.Item(2, 11).Value = 1.007221551487641
$ws.Cells.Item(2, 12).Value = 1.006920690104309
$ws.Cells.Item(2, 13).Value = 1.006875393711417
$ws.Cells.Item(2, 14).Value = 1.022734980948664
$ws.Cells.Item(2, 15).Value = 1.021430049851221
$ws.Cells.Item(2, 16).Value = 1.02067649724301
$ws.Cells.Item(2, 17).Value = 1.007446506244054
$ws.Cells.Item(2, 18).Value = 1.021413915781757
$ws.Cells.Item(2, 19).Value = 1.008002208741513
$ws.Cells.Item(2, 20).Value = 1.006859489541038
$ws.Cells.Item(2, 21).Value = 1.008748309550161
$ws.Cells.Item(3, 3).Value = 1.022654000256572
$ws.Cells.Item(3, 4).Value = 1.017826244806614
$ws.Cells.Item(3, 5).Value = 1.009760920452753
$ws.Cells.Item(3, 6).Value = 1.009336673118919
$ws.Cells.Item(3, 7).Value = 1.009047031252258
$ws.Cells.Item(3, 8).Value = 1.008712747020375
$ws.Cells.Item(3, 9).Value = 1.008482582963286
$ws.Cells.Item(3, 10).Value = 1.008516235247874
$ws.Cells.Item(3, 11).Value = 1.008285332870445
$ws.Cells.Item(3, 12).Value = 1.007997198739165
$ws.Cells.Item(3, 13).Value = 1.007953904000022
$ws.Cells.Item(3, 14).Value = 1.023023459945867
$ws.Cells.Item(3, 15).Value = 1.021768915316444
$ws.Cells.Item(3, 16).Value = 1.021044215156406
$ws.Cells.Item(3, 17).Value = 1.008500305159468
$ws.Cells.Item(3, 18).Value = 1.021752775894406
$ws.Cells.Item(3, 19).Value = 1.009031092779627
$ws.Cells.Item(3, 20).Value = 1.007937982793957
$ws.Cells.Item(3, 21).Value = 1.009744970703834
$ws.Cells.Item(4, 3).Value = 1.023030459902508
$ws.Cells.Item(4, 4).Value = 1.018457524838795
$ws.Cells.Item(4, 5).Value = 1.010787487972367
$ws.Cells.Item(4, 6).Value = 1.010382957995344
$ws.Cells.Item(4, 7).Value = 1.010106839194546
$ws.Cells.Item(4, 8).Value = 1.009788545020354
$ws.Cells.Item(4, 9).Value = 1.009571488756242
$ws.Cells.Item(4, 10).Value = 1.00960174003783
$ws.Cells.Item(4, 11).Value = 1.009381117377024
$ws.Cells.Item(4, 12).Value = 1.009106111240664
$ws.Cells.Item(4, 13).Value = 1.009064881190704
$ws.Cells.Item(4, 14).Value = 1.023320448473356
$ws.Cells.Item(4, 15).Value = 1.022117955738541
$ws.Cells.Item(4, 16).Value = 1.021423061652893
$ws.Cells.Item(4, 17).Value = 1.009585792803257
$ws.Cells.Item(4, 18).Value = 1.022101810803212
$ws.Cells.Item(4, 19).Value = 1.010090883981644
$ws.Cells.Item(4, 20).Value = 1.009048942436122
$ws.Cells.Item(4, 21).Value = 1.01077152200823
$ws.Cells.Item(5, 3).Value = 1.023032344977129
$ws.Cells.Item(5, 4).Value = 1.018460687785299
$ws.Cells.Item(5, 5).Value = 1.010792633532484
$ws.Cells.Item(5, 6).Value = 1.010388202463306
$ws.Cells.Item(5, 7).Value = 1.010112151498272
$ws.Cells.Item(5, 8).Value = 1.009793937534579
$ws.Cells.Item(5, 9).Value = 1.009576947010703
$ws.Cells.Item(5, 10).Value = 1.009607181232206
$ws.Cells.Item(5, 11).Value = 1.009386610135802
$ws.Cells.Item(5, 12).Value = 1.009111669850853
$ws.Cells.Item(5, 13).Value = 1.009070450157617
$ws.Cells.Item(5, 14).Value = 1.023321936726942
$ws.Cells.Item(5, 15).Value = 1.022119705290602
$ws.Cells.Item(5, 16).Value = 1.021424960830191
$ws.Cells.Item(5, 17).Value = 1.009591233911687
$ws.Cells.Item(5, 18).Value = 1.022103560327637
$ws.Cells.Item(5, 19).Value = 1.010096196201459
$ws.Cells.Item(5, 20).Value = 1.00905451131507
$ws.Cells.Item(5, 21).Value = 1.01077666748707
$ws.Cells.Item(6, 3).Value = 1.023273096800445
$ws.Cells.Item(6, 4).Value = 1.018864800888
$ws.Cells.Item(6, 5).Value = 1.01145022836299
$ws.Cells.Item(6, 6).Value = 1.011058443968311
$ws.Cells.Item(6, 7).Value = 1.010791066721088
$ws.Cells.Item(6, 8).Value = 1.010483108715772
$ws.Cells.Item(6, 9).Value = 1.010274522940927
$ws.Cells.Item(6, 10).Value = 1.010302575830808
$ws.Cells.Item(6, 11).Value = 1.010088597794501
$ws.Cells.Item(6, 12).Value = 1.00982207728906
$ws.Cells.Item(6, 13).Value = 1.009782181811265
$ws.Cells.Item(6, 14).Value = 1.023512101775912
$ws.Cells.Item(6, 15).Value = 1.02234329669413
$ws.Cells.Item(6, 16).Value = 1.021667692889974
$ws.Cells.Item(6, 17).Value = 1.010286617526135
$ws.Cells.Item(6, 18).Value = 1.022327148199412
$ws.Cells.Item(6, 19).Value = 1.010775100700422
$ws.Cells.Item(6, 20).Value = 1.009766231726512
$ws.Cells.Item(6, 21).Value = 1.011434251930491
$ws.Cells.Item(7, 3).Value = 1.023085262169533
$ws.Cells.Item(7, 4).Value = 1.018549485018173
$ws.Cells.Item(7, 5).Value = 1.010937099903137
$ws.Cells.Item(7, 6).Value = 1.010535446088573
$ws.Cells.Item(7, 7).Value = 1.010261299899466
$ws.Cells.Item(7, 8).Value = 1.009945338182877
$ws.Cells.Item(7, 9).Value = 1.00973019354211
$ws.Cells.Item(7, 10).Value = 1.009759948731113
$ws.Cells.Item(7, 11).Value = 1.009540825516052
$ws.Cells.Item(7, 12).Value = 1.009267734277196
$ws.Cells.Item(7, 13).Value = 1.009226805391169
$ws.Cells.Item(7, 14).Value = 1.023363719215656
$ws.Cells.Item(7, 15).Value = 1.02216882564266
$ws.Cells.Item(7, 16).Value = 1.021478283003358
$ws.Cells.Item(7, 17).Value = 1.009743998997543
$ws.Cells.Item(7, 18).Value = 1.022152679903812
$ws.Cells.Item(7, 19).Value = 1.010245342246769
$ws.Cells.Item(7, 20).Value = 1.009210864078903
$ws.Cells.Item(7, 21).Value = 1.010921131575795
$ws.Cells.Item(8, 3).Value = 1.022630780272726
$ws.Cells.Item(8, 4).Value = 1.01778733211316
$ws.Cells.Item(8, 5).Value = 1.00969766887615
$ws.Cells.Item(8, 6).Value = 1.009272207662286
$ws.Cells.Item(8, 7).Value = 1.008981733264197
$ws.Cells.Item(8, 8).Value = 1.008646464617772
$ws.Cells.Item(8, 9).Value = 1.008415493433572
$ws.Cells.Item(8, 10).Value = 1.008449355101644
$ws.Cells.Item(8, 11).Value = 1.008217819842967
$ws.Cells.Item(8, 12).Value = 1.007928877468999
$ws.Cells.Item(8, 13).Value = 1.007885455614795
$ws.Cells.Item(8, 14).Value = 1.023005156338313
$ws.Cells.Item(8, 15).Value = 1.021747409679758
$ws.Cells.Item(8, 16).Value = 1.021020875965341
$ws.Cells.Item(8, 17).Value = 1.008433426069647
$ws.Cells.Item(8, 18).Value = 1.021731270597415
$ws.Cells.Item(8, 19).Value = 1.008965795822985
$ws.Cells.Item(8, 20).Value = 1.007869535489912
$ws.Cells.Item(8, 21).Value = 1.009681720126326
$ws.Cells.Item(9, 3).Value = 1.021363395167066
$ws.Cells.Item(9, 4).Value = 1.015644861136213
$ws.Cells.Item(9, 5).Value = 1.006199090528654
$ws.Cells.Item(9, 6).Value = 1.005707880068399
$ws.Cells.Item(9, 7).Value = 1.005372312444564
$ws.Cells.Item(9, 8).Value = 1.00498372247808
$ws.Cells.Item(9, 9).Value = 1.004698221098066
$ws.Cells.Item(9, 10).Value = 1.004743428567433
$ws.Cells.Item(9, 11).Value = 1.004474964022374
$ws.Cells.Item(9, 12).Value = 1.004135879433539
$ws.Cells.Item(9, 13).Value = 1.004082849790051
$ws.Cells.Item(9, 14).Value = 1.022010969472218
$ws.Cells.Item(9, 15).Value = 1.020580323946099
$ws.Cells.Item(9, 16).Value = 1.019754784223468
$ws.Cells.Item(9, 17).Value = 1.004727558072657
$ws.Cells.Item(9, 18).Value = 1.02056420329854
$ws.Cells.Item(9, 19).Value = 1.005356432016209
$ws.Cells.Item(9, 20).Value = 1.003940932901547
$ws.Cells.Item(9, 21).Value = 1.006183197040868
$ws.Cells.Item(10, 3).Value = 1.016814234525497
$ws.Cells.Item(10, 4).Value = 1.008122392769345
$ws.Cells.Item(10, 5).Value = 0.9940787794780812
$ws.Cells.Item(10, 6).Value = 0.993356963544071
$ws.Cells.Item(10, 7).Value = 0.9928632470425405
$ws.Cells.Item(10, 8).Value = 0.992287637653297
$ws.Cells.Item(10, 9).Value = 0.9918590550176571
$ws.Cells.Item(10, 10).Value = 0.9919440492234796
$ws.Cells.Item(10, 11).Value = 0.9915579309088496
$ws.Cells.Item(10, 12).Value = 0.9910715578341394
$ws.Cells.Item(10, 13).Value = 0.9909970864997731
$ws.Cells.Item(10, 14).Value = 1.018471851288541
$ws.Cells.Item(10, 15).Value = 1.016441188574724
$ws.Cells.Item(10, 16).Value = 1.015272061620074
$ws.Cells.Item(10, 17).Value = 0.9919283809021914
$ws.Cells.Item(10, 18).Value = 1.016425133307166
$ws.Cells.Item(10, 19).Value = 0.992847564201999
$ws.Cells.Item(10, 20).Value = 0.9909814331363003
$ws.Cells.Item(10, 21).Value = 0.9940630774375127
$ws.Cells.Item(11, 3).Value = 1.012579261596185
$ws.Cells.Item(11, 4).Value = 1.001129543381243
$ws.Cells.Item(11, 5).Value = 0.9828381145814358
$ws.Cells.Item(11, 6).Value = 0.9819042122236248
$ws.Cells.Item(11, 7).Value = 0.9812638691735586
$ws.Cells.Item(11, 8).Value = 0.9805204341535008
$ws.Cells.Item(11, 9).Value = 0.9799631021533486
$ws.Cells.Item(11, 10).Value = 0.9800838190276002
$ws.Cells.Item(11, 11).Value = 0.9795896827800997
$ws.Cells.Item(11, 12).Value = 0.9789653319239022
$ws.Cells.Item(11, 13).Value = 0.9788691622574641
$ws.Cells.Item(11, 14).Value = 1.015240276294989
$ws.Cells.Item(11, 15).Value = 1.012681064882454
$ws.Cells.Item(11, 16).Value = 1.011209271566258
$ws.Cells.Item(11, 17).Value = 0.9800683380454035
$ws.Cells.Item(11, 18).Value = 1.012665069008192
$ws.Cells.Item(11, 19).Value = 0.9811193715294424
$ws.Cells.Item(11, 20).Value = 0.978853700461463
$ws.Cells.Item(11, 21).Value = 0.9826937993875077
$ws.Cells.Item(12, 3).Value = 1.005297856383255
$ws.Cells.Item(12, 4).Value = 0.9893516215900886
$ws.Cells.Item(12, 5).Value = 0.9641629229935793
$ws.Cells.Item(12, 6).Value = 0.9628915279497178
$ws.Cells.Item(12, 7).Value = 0.9620209654525572
$ws.Cells.Item(12, 8).Value = 0.9609996113219076
$ws.Cells.Item(12, 9).Value = 0.96020144380735
$ws.Cells.Item(12, 10).Value = 0.9603813386603203
$ws.Cells.Item(12, 11).Value = 0.9597081464970401
$ws.Cells.Item(12, 12).Value = 0.9588550439133635
$ws.Cells.Item(12, 13).Value = 0.9587229041100579
$ws.Cells.Item(12, 14).Value = 1.009783371471103
$ws.Cells.Item(12, 15).Value = 1.006368473770559
$ws.Cells.Item(12, 16).Value = 1.004406674851279
$ws.Cells.Item(12, 17).Value = 0.9603661688900241
$ws.Cells.Item(12, 18).Value = 1.006352577607268
$ws.Cells.Item(12, 19).Value = 0.9620057697834216
$ws.Cells.Item(12, 20).Value = 0.95870776053568
$ws.Cells.Item(12, 21).Value = 0.96401639955618
$ws.Cells.Item(13, 3).Value = 0.9992007685579863
$ws.Cells.Item(13, 4).Value = 0.9795815468474661
$ws.Cells.Item(13, 5).Value = 0.948786756158456
$ws.Cells.Item(13, 6).Value = 0.9472342629599226
$ws.Cells.Item(13, 7).Value = 0.9461697523793638
$ws.Cells.Item(13, 8).Value = 0.9449247364735041
$ws.Cells.Item(13, 9).Value = 0.9439498246874815
$ws.Cells.Item(13, 10).Value = 0.9441773662652847
$ws.Cells.Item(13, 11).Value = 0.9433599359561322
$ws.Cells.Item(13, 12).Value = 0.9423225294464589
$ws.Cells.Item(13, 13).Value = 0.9421614112710135
$ws.Cells.Item(13, 14).Value = 1.005323835672268
$ws.Cells.Item(13, 15).Value = 1.00123990400822
$ws.Cells.Item(13, 16).Value = 0.9988949779369247
$ws.Cells.Item(13, 17).Value = 0.9441624524459624
$ws.Cells.Item(13, 18).Value = 1.00122408885361
$ws.Cells.Item(13, 19).Value = 0.9460210073772651
$ws.Cells.Item(13, 20).Value = 0.9421465292948485
$ws.Cells.Item(13, 21).Value = 0.9487717695311947
$ws.Cells.Item(14, 3).Value = 0.9986215517192463
$ws.Cells.Item(14, 4).Value = 0.9786578776209675
$ws.Cells.Item(14, 5).Value = 0.9473385853787518
$ws.Cells.Item(14, 6).Value = 0.9457603423155501
$ws.Cells.Item(14, 7).Value = 0.9446781772519762
$ws.Cells.Item(14, 8).Value = 0.943412280980801
$ws.Cells.Item(14, 9).Value = 0.9424198641885545
$ws.Cells.Item(14, 10).Value = 0.9426518552720584
$ws.Cells.Item(14, 11).Value = 0.9418209512915393
$ws.Cells.Item(14, 12).Value = 0.9407663268951646
$ws.Cells.Item(14, 13).Value = 0.9406025016556544
$ws.Cells.Item(14, 14).Value = 1.004905102963169
$ws.Cells.Item(14, 15).Value = 1.000759622930117
$ws.Cells.Item(14, 16).Value = 0.9983794504093206
$ws.Cells.Item(14, 17).Value = 0.9426369655490513
$ws.Cells.Item(14, 18).Value = 1.000743815361821
$ws.Cells.Item(14, 19).Value = 0.944529243781845
$ws.Cells.Item(14, 20).Value = 0.9405876443033552
$ws.Cells.Item(14, 21).Value = 0.9473236216261729
$ws.Cells.Item(15, 3).Value = 0.9990932648265023
$ws.Cells.Item(15, 4).Value = 0.9794100578919419
$ws.Cells.Item(15, 5).Value = 0.9485178203192088
$ws.Cells.Item(15, 6).Value = 0.9469605423227835
$ws.Cells.Item(15, 7).Value = 0.9458927512012324
$ws.Cells.Item(15, 8).Value = 0.9446438553476416
$ws.Cells.Item(15, 9).Value = 0.9436656915825906
$ws.Cells.Item(15, 10).Value = 0.9438940599440638
$ws.Cells.Item(15, 11).Value = 0.9430741260810747
$ws.Cells.Item(15, 12).Value = 0.9420335203224284
$ws.Cells.Item(15, 13).Value = 0.9418718991467286
$ws.Cells.Item(15, 14).Value = 1.005246053913137
$ws.Cells.Item(15, 15).Value = 1.001150673458975
$ws.Cells.Item(15, 16).Value = 0.9987991910979196
$ws.Cells.Item(15, 17).Value = 0.9438791505997262
$ws.Cells.Item(15, 18).Value = 1.001134859713812
$ws.Cells.Item(15, 19).Value = 0.945743971249475
$ws.Cells.Item(15, 20).Value = 0.9418570217435728
$ws.Cells.Item(15, 21).Value = 0.9485028379399425
$ws.Cells.Item(16, 3).Value = 0.9991154903043679
$ws.Cells.Item(16, 4).Value = 0.9794698238863061
$ws.Cells.Item(16, 5).Value = 0.9486338037481729
$ws.Cells.Item(16, 6).Value = 0.9470826784383699
$ws.Cells.Item(16, 7).Value = 0.9460203364731335
$ws.Cells.Item(16, 8).Value = 0.9447720817729454
$ws.Cells.Item(16, 9).Value = 0.9437841314117045
$ws.Cells.Item(16, 10).Value = 0.9440123636415927
$ws.Cells.Item(16, 11).Value = 0.9431928329934556
$ws.Cells.Item(16, 12).Value = 0.9421527429998072
$ws.Cells.Item(16, 13).Value = 0.9419912030542555
$ws.Cells.Item(16, 14).Value = 1.005256294824117
$ws.Cells.Item(16, 15).Value = 1.001162421326962
$ws.Cells.Item(16, 16).Value = 0.9988118019433887
$ws.Cells.Item(16, 17).Value = 0.943997452428581
$ws.Cells.Item(16, 18).Value = 1.001146607396236
$ws.Cells.Item(16, 19).Value = 0.9460053935431979
$ws.Cells.Item(16, 20).Value = 0.9419763237666264
$ws.Cells.Item(16, 21).Value = 0.9486188195368824
$ws.Cells.Item(17, 3).Value = 0.9991685450904111
$ws.Cells.Item(17, 4).Value = 0.9795298811143368
$ws.Cells.Item(17, 5).Value = 0.9487054610845316
$ws.Cells.Item(17, 6).Value = 0.947156925170781
$ws.Cells.Item(17, 7).Value = 0.9460963582884252
$ws.Cells.Item(17, 8).Value = 0.9448502033189333
$ws.Cells.Item(17, 9).Value = 0.9438523826731062
$ws.Cells.Item(17, 10).Value = 0.9440801740668923
$ws.Cells.Item(17, 11).Value = 0.9432591592803155
$ws.Cells.Item(17, 12).Value = 0.9422139971591037
$ws.Cells.Item(17, 13).Value = 0.9420498245984484
$ws.Cells.Item(17, 14).Value = 1.005300749529178
$ws.Cells.Item(17, 15).Value = 1.0012134190303
$ws.Cells.Item(17, 16).Value = 0.9988665465748486
$ws.Cells.Item(17, 17).Value = 0.9440652617827765
$ws.Cells.Item(17, 18).Value = 1.001197604294035
$ws.Cells.Item(17, 19).Value = 0.946081414157531
$ws.Cells.Item(17, 20).Value = 0.941900557379894
$ws.Cells.Item(17, 21).Value = 0.9486904757413724
$ws.Cells.Item(18, 3).Value = 1.005036260194141
$ws.Cells.Item(18, 4).Value = 0.9889541967628399
$ws.Cells.Item(18, 5).Value = 0.9635569414156528
$ws.Cells.Item(18, 6).Value = 0.9622733393025414
$ws.Cells.Item(18, 7).Value = 0.9613944074160036
$ws.Cells.Item(18, 8).Value = 0.9603631525905185
$ws.Cells.Item(18, 9).Value = 0.9595567316133728
$ws.Cells.Item(18, 10).Value = 0.9597387327633737
$ws.Cells.Item(18, 11).Value = 0.9590591572655405
$ws.Cells.Item(18, 12).Value = 0.9581978993226578
$ws.Cells.Item(18, 13).Value = 0.9580644773226147
$ws.Cells.Item(18, 14).Value = 1.009585311052824
$ws.Cells.Item(18, 15).Value = 1.006140152316456
$ws.Cells.Item(18, 16).Value = 1.004161024335809
$ws.Cells.Item(18, 17).Value = 0.9597235731434036
$ws.Cells.Item(18, 18).Value = 1.006124259759633
$ws.Cells.Item(18, 19).Value = 0.9613792216437091
$ws.Cells.Item(18, 20).Value = 0.9580493441484635
$ws.Cells.Item(18, 21).Value = 0.9635417214849021
$ws.Cells.Item(19, 3).Value = 1.014644825834384
$ws.Cells.Item(19, 4).Value = 1.0044859306647
$ws.Cells.Item(19, 5).Value = 0.9881818418478231
$ws.Cells.Item(19, 6).Value = 0.9873535109780535
$ws.Cells.Item(19, 7).Value = 0.9867867509051887
$ws.Cells.Item(19, 8).Value = 0.9861247533970452
$ws.Cells.Item(19, 9).Value = 0.9856027626734759
$ws.Cells.Item(19, 10).Value = 0.9857061536823657
$ws.Cells.Item(19, 11).Value = 0.985259026330925
$ws.Cells.Item(19, 12).Value = 0.9846886414928995
$ws.Cells.Item(19, 13).Value = 0.9845974429970322
$ws.Cells.Item(19, 14).Value = 1.016814662975463
$ws.Cells.Item(19, 15).Value = 1.014470401983127
$ws.Cells.Item(19, 16).Value = 1.013144948709212
$ws.Cells.Item(19, 17).Value = 0.9856905838921906
$ws.Cells.Item(19, 18).Value = 1.014329615859043
$ws.Cells.Item(19, 19).Value = 0.9867711640463649
$ws.Cells.Item(19, 20).Value = 0.9843243399380063
$ws.Cells.Item(19, 21).Value = 0.9880381410241772
$ws.Cells.Item(20, 3).Value = 1.017021268077087
$ws.Cells.Item(20, 4).Value = 1.008417848860809
$ws.Cells.Item(20, 5).Value = 0.9945116483922374
$ws.Cells.Item(20, 6).Value = 0.9938004062075536
$ws.Cells.Item(20, 7).Value = 0.9933139407192726
$ws.Cells.Item(20, 8).Value = 0.9927469055560827
$ws.Cells.Item(20, 9).Value = 0.9923144417915316
$ws.Cells.Item(20, 10).Value = 0.9923976180917748
$ws.Cells.Item(20, 11).Value = 0.9920143350349426
$ws.Cells.Item(20, 12).Value = 0.9915285654837306
$ws.Cells.Item(20, 13).Value = 0.9914524451902142
$ws.Cells.Item(20, 14).Value = 1.018637431352431
$ws.Cells.Item(20, 15).Value = 1.016634328864176
$ws.Cells.Item(20, 16).Value = 1.015480983800037
$ws.Cells.Item(20, 17).Value = 0.9923819426061079
$ws.Cells.Item(20, 18).Value = 1.016618270545856
$ws.Cells.Item(20, 19).Value = 0.9932982507597676
$ws.Cells.Item(20, 20).Value = 0.991309116654614
$ws.Cells.Item(20, 21).Value = 0.9943686655449495
$ws.Cells.Item(21, 3).Value = 1.019493935132369
$ws.Cells.Item(21, 4).Value = 1.01253269180433
$ws.Cells.Item(21, 5).Value = 1.00116330884761
$ws.Cells.Item(21, 6).Value = 1.000571246704534
$ws.Cells.Item(21, 7).Value = 1.00016534495993
$ws.Cells.Item(21, 8).Value = 0.9996990672186851
$ws.Cells.Item(21, 9).Value = 0.999370229859364
$ws.Cells.Item(21, 10).Value = 0.9994319618082639
$ws.Cells.Item(21, 11).Value = 0.999116185211681
$ws.Cells.Item(21, 12).Value = 0.9987196546569652
$ws.Cells.Item(21, 13).Value = 0.9986593122481833
$ws.Cells.Item(21, 14).Value = 1.020550319408322
$ws.Cells.Item(21, 15).Value = 1.01886919354919
$ws.Cells.Item(21, 16).Value = 1.017900222117251
$ws.Cells.Item(21, 17).Value = 0.999416175211131
$ws.Cells.Item(21, 18).Value = 1.01885309992991
$ws.Cells.Item(21, 19).Value = 1.000022994602922
$ws.Cells.Item(21, 20).Value = 0.9986435378554902
$ws.Cells.Item(21, 21).Value = 1.001147494902864
$ws.Cells.Item(22, 3).Value = 1.020387907091106
$ws.Cells.Item(22, 4).Value = 1.014018869207036
$ws.Cells.Item(22, 5).Value = 1.003565790772829
$ws.Cells.Item(22, 6).Value = 1.003019350018149
$ws.Cells.Item(22, 7).Value = 1.002644738489521
$ws.Cells.Item(22, 8).Value = 1.002215448490336
$ws.Cells.Item(22, 9).Value = 1.001916990819399
$ws.Cells.Item(22, 10).Value = 1.001970853808667
$ws.Cells.Item(22, 11).Value = 1.001678866124234
$ws.Cells.Item(22, 12).Value = 1.00131271818597
$ws.Cells.Item(22, 13).Value = 1.001257153995579
$ws.Cells.Item(22, 14).Value = 1.021247445476572
$ws.Cells.Item(22, 15).Value = 1.019685351916267
$ws.Cells.Item(22, 16).Value = 1.018784541653293
$ws.Cells.Item(22, 17).Value = 1.001955027108289
$ws.Cells.Item(22, 18).Value = 1.019669245405301
$ws.Cells.Item(22, 19).Value = 1.002502662927118
$ws.Cells.Item(22, 20).Value = 1.001241338568496
$ws.Cells.Item(22, 21).Value = 1.003549938879513
$ws.Cells.Item(23, 3).Value = 1.020691703352726
$ws.Cells.Item(23, 4).Value = 1.01454752307477
$ws.Cells.Item(23, 5).Value = 1.004442181352812
$ws.Cells.Item(23, 6).Value = 1.003916258241095
$ws.Cells.Item(23, 7).Value = 1.00355688396569
$ws.Cells.Item(23, 8).Value = 1.003140143366839
$ws.Cells.Item(23, 9).Value = 1.002842309368001
$ws.Cells.Item(23, 10).Value = 1.002893501435303
$ws.Cells.Item(23, 11).Value = 1.002609580027543
$ws.Cells.Item(23, 12).Value = 1.002253734394235
$ws.Cells.Item(23, 13).Value = 1.002199790564773
$ws.Cells.Item(23, 14).Value = 1.021482550812523
$ws.Cells.Item(23, 15).Value = 1.019960811448615
$ws.Cells.Item(23, 16).Value = 1.019083109065594
$ws.Cells.Item(23, 17).Value = 1.00287766016118
$ws.Cells.Item(23, 18).Value = 1.019944700586609
$ws.Cells.Item(23, 19).Value = 1.00354103221437
$ws.Cells.Item(23, 20).Value = 1.002183960248208
$ws.Cells.Item(23, 21).Value = 1.004426315616407
$ws.Cells.Item(24, 3).Value = 1.02143047738283
$ws.Cells.Item(24, 4).Value = 1.015779633355027
$ws.Cells.Item(24, 5).Value = 1.006438391812122
$ws.Cells.Item(24, 6).Value = 1.005950535740439
$ws.Cells.Item(24, 7).Value = 1.005617268730074
$ws.Cells.Item(24, 8).Value = 1.00523139969446
$ws.Cells.Item(24, 9).Value = 1.004958916346667
$ws.Cells.Item(24, 10).Value = 1.005003541199729
$ws.Cells.Item(24, 11).Value = 1.004739471894131
$ws.Cells.Item(24, 12).Value = 1.004408980128753
$ws.Cells.Item(24, 13).Value = 1.004359023710278
$ws.Cells.Item(24, 14).Value = 1.022061284486923
$ws.Cells.Item(24, 15).Value = 1.020639341486578
$ws.Cells.Item(24, 16).Value = 1.019818785046902
$ws.Cells.Item(24, 17).Value = 1.004987666596326
$ws.Cells.Item(24, 18).Value = 1.020623219906803
$ws.Cells.Item(24, 19).Value = 1.005601384432494
$ws.Cells.Item(24, 20).Value = 1.004343159287396
$ws.Cells.Item(24, 21).Value = 1.006422494544436
$ws.Cells.Item(25, 3).Value = 1.021759745374283
$ws.Cells.Item(25, 4).Value = 1.016329653648021
$ws.Cells.Item(25, 5).Value = 1.00733047989036
$ws.Cells.Item(25, 6).Value = 1.00685967173793
$ws.Cells.Item(25, 7).Value = 1.006538096703128
$ws.Cells.Item(25, 8).Value = 1.006166053090969
$ws.Cells.Item(25, 9).Value = 1.005904916820103
$ws.Cells.Item(25, 10).Value = 1.005946600764998
$ws.Cells.Item(25, 11).Value = 1.005691421245162
$ws.Cells.Item(25, 12).Value = 1.005372282800053
$ws.Cells.Item(25, 13).Value = 1.005324111933676
$ws.Cells.Item(25, 14).Value = 1.022319763117535
$ws.Cells.Item(25, 15).Value = 1.0209426063822
$ws.Cells.Item(25, 16).Value = 1.020147695657129
$ws.Cells.Item(25, 17).Value = 1.005930711265432
$ws.Cells.Item(25, 18).Value = 1.020926480012183
$ws.Cells.Item(25, 19).Value = 1.006522197860547
$ws.Cells.Item(25, 20).Value = 1.005308232266676
$ws.Cells.Item(25, 21).Value = 1.007314568531635
